# Auto-generated: update Leve profit-tracking values per scheduled market-price refresh
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3437.3572
$ws.Range("I64").Value = 3394
$ws.Range("J64").Value = 3495.1667
$ws.Range("K64").Value = 3394
$ws.Range("L64").Value = 3495.1667
$ws.Range("M64").Value = -3146
$ws.Range("N64").Value = -3991.1667

$ws.Range("H67").Value = 3437.3572
$ws.Range("I67").Value = 3394
$ws.Range("J67").Value = 3495.1667
$ws.Range("K67").Value = 3394
$ws.Range("L67").Value = 3495.1667
$ws.Range("M67").Value = -2536
$ws.Range("N67").Value = -5211.1667

$ws.Range("H76").Value = 3090406
$ws.Range("I76").Value = 4118879.5
$ws.Range("J76").Value = 4985
$ws.Range("K76").Value = 4118879.5
$ws.Range("L76").Value = 4985
$ws.Range("M76").Value = -4118564.5
$ws.Range("N76").Value = -5615

$ws.Range("H79").Value = 3090406
$ws.Range("I79").Value = 4118879.5
$ws.Range("J79").Value = 4985
$ws.Range("K79").Value = 4118879.5
$ws.Range("L79").Value = 4985
$ws.Range("M79").Value = -4117787.5
$ws.Range("N79").Value = -7169

$ws.Range("H94").Value = 2431.0908
$ws.Range("I94").Value = 2175.5
$ws.Range("J94").Value = 4987
$ws.Range("K94").Value = 2175.5
$ws.Range("L94").Value = 4987
$ws.Range("M94").Value = -1724.5
$ws.Range("N94").Value = -5889

$ws.Range("H101").Value = 532.55554
$ws.Range("I101").Value = 451
$ws.Range("J101").Value = 1185
$ws.Range("K101").Value = 1353
$ws.Range("L101").Value = 3555
$ws.Range("M101").Value = 269
$ws.Range("N101").Value = -6799

$ws.Range("H113").Value = 2831.2258
$ws.Range("I113").Value = 2459.2
$ws.Range("J113").Value = 3180
$ws.Range("K113").Value = 2459.2
$ws.Range("L113").Value = 3180
$ws.Range("M113").Value = 794.8000000000002
$ws.Range("N113").Value = -9688

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3800.25
$ws.Range("I63").Value = 2164
$ws.Range("J63").Value = 7400
$ws.Range("K63").Value = 2164
$ws.Range("L63").Value = 7400
$ws.Range("M63").Value = -1478
$ws.Range("N63").Value = -8772

$ws.Range("H66").Value = 3800.25
$ws.Range("I66").Value = 2164
$ws.Range("J66").Value = 7400
$ws.Range("K66").Value = 10820
$ws.Range("L66").Value = 37000
$ws.Range("M66").Value = -7388
$ws.Range("N66").Value = -43864

$ws.Range("H88").Value = 5024.2104
$ws.Range("I88").Value = 1766.6666
$ws.Range("J88").Value = 5635
$ws.Range("K88").Value = 1766.6666
$ws.Range("L88").Value = 5635
$ws.Range("M88").Value = -1360.6666
$ws.Range("N88").Value = -6447

$ws.Range("H91").Value = 5024.2104
$ws.Range("I91").Value = 1766.6666
$ws.Range("J91").Value = 5635
$ws.Range("K91").Value = 1766.6666
$ws.Range("L91").Value = 5635
$ws.Range("M91").Value = -362.6666
$ws.Range("N91").Value = -8443

$ws.Range("H131").Value = 62692
$ws.Range("J131").Value = 62692
$ws.Range("L131").Value = 62692
$ws.Range("N131").Value = -72772

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4505.2563
$ws.Range("I31").Value = 2975.6584
$ws.Range("J31").Value = 6200.2163
$ws.Range("K31").Value = 2975.6584
$ws.Range("L31").Value = 6200.2163
$ws.Range("M31").Value = -2680.6584
$ws.Range("N31").Value = -6790.2163

$ws.Range("H34").Value = 4505.2563
$ws.Range("I34").Value = 2975.6584
$ws.Range("J34").Value = 6200.2163
$ws.Range("K34").Value = 2975.6584
$ws.Range("L34").Value = 6200.2163
$ws.Range("M34").Value = -2773.6584
$ws.Range("N34").Value = -6604.2163

$ws.Range("H62").Value = 3478.6316
$ws.Range("I62").Value = 3708.8462
$ws.Range("J62").Value = 2979.8333
$ws.Range("K62").Value = 3708.8462
$ws.Range("L62").Value = 2979.8333
$ws.Range("M62").Value = -3084.8462
$ws.Range("N62").Value = -4227.8333

$ws.Range("H65").Value = 3478.6316
$ws.Range("I65").Value = 3708.8462
$ws.Range("J65").Value = 2979.8333
$ws.Range("K65").Value = 18544.231
$ws.Range("L65").Value = 14899.1665
$ws.Range("M65").Value = -15424.231
$ws.Range("N65").Value = -21139.1665

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 9555.556

$ws.Range("H90").Value = 9555.556

$ws.Range("H100").Value = 5000
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 15000
$ws.Range("N100").Value = -16622

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2854.5
$ws.Range("I80").Value = 2693.6
$ws.Range("J80").Value = 2969.4285
$ws.Range("K80").Value = 2693.6
$ws.Range("L80").Value = 2969.4285
$ws.Range("M80").Value = -1695.6
$ws.Range("N80").Value = -4965.4285

$ws.Range("H83").Value = 2854.5
$ws.Range("I83").Value = 2693.6
$ws.Range("J83").Value = 2969.4285
$ws.Range("K83").Value = 13468
$ws.Range("L83").Value = 14847.1425
$ws.Range("M83").Value = -8476
$ws.Range("N83").Value = -24831.1425

$ws.Range("H102").Value = 1595.9333
$ws.Range("I102").Value = 1516.2916
$ws.Range("J102").Value = 1914.5
$ws.Range("K102").Value = 1516.2916
$ws.Range("L102").Value = 1914.5
$ws.Range("M102").Value = 105.7084
$ws.Range("N102").Value = -5158.5

$ws.Range("H113").Value = 71430540
$ws.Range("I113").Value = 1962.5555
$ws.Range("J113").Value = 200001980
$ws.Range("K113").Value = 1962.5555
$ws.Range("L113").Value = 200001980
$ws.Range("M113").Value = 207.4445000000001
$ws.Range("N113").Value = -200006320

$ws.Range("H126").Value = 8773315
$ws.Range("I126").Value = 1600.2222
$ws.Range("J126").Value = 16667858
$ws.Range("K126").Value = 4800.6666
$ws.Range("L126").Value = 50003574
$ws.Range("M126").Value = -2330.6666
$ws.Range("N126").Value = -50008514

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()

$ws.Range("H102").Value = 36000
$ws.Range("J102").Value = 36000
$ws.Range("L102").Value = 36000
$ws.Range("N102").Value = -42490

$ws.Range("H122").Value = 2472.7942
$ws.Range("I122").Value = 1976.9048
$ws.Range("J122").Value = 3273.8462
$ws.Range("K122").Value = 5930.7144
$ws.Range("L122").Value = 9821.5386
$ws.Range("M122").Value = -3480.7144
$ws.Range("N122").Value = -14721.5386

$ws.Range("H136").Value = 3720.9333
$ws.Range("I136").Value = 2350.7827
$ws.Range("K136").Value = 7052.348100000001
$ws.Range("M136").Value = -4502.348100000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 592.12067
$ws.Range("I136").Value = 513.04254
$ws.Range("J136").Value = 930
$ws.Range("K136").Value = 1539.12762
$ws.Range("L136").Value = 2790
$ws.Range("M136").Value = 1010.87238
$ws.Range("N136").Value = -7890

